# Update the "For plotting" sheet's D2:G16 range (Lower, Upper, Bias_avg, Bias_se)
# with refreshed bootstrap estimates from the re-run clustering analysis
# (clustering without repeated families, 0 bias and less singular fits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("For plotting")

$values = @(
    @(0.299884212478717,0.405331310068447,-0.0028946039655305,0.0268997697932983),
    @(0.195218193880211,0.288758868116803,-0.00184157646055815,0.0238624168970899),
    @(0.383058767492031,0.482958430273897,-0.00243416950137019,0.0254846078525168),
    @(0.07016645861601,0.177112270642956,-0.00430274190057132,0.0272820949048331),
    @(0.0196314615567271,0.129091904749199,-0.00451263285412443,0.0279235824470593),
    @(0.137878614580662,0.258604411735044,-0.000770376762267614,0.0307973972332607),
    @(0.0978387504233838,0.236851713375691,0.000278136555685943,0.035462490549058),
    @(0.0907822240865411,0.221605276366246,-0.00483062708327292,0.0333732276223736),
    @(0.155837072726262,0.288144652159785,-0.00481561623999976,0.0337519335289598),
    @(0.0913416363818949,0.370815978486343,-0.00864538559713109,0.071294475026645),
    @(0.0700554468988439,0.336586106552234,-0.00927909424664185,0.0679925152177015),
    @(0.102374803560793,0.374311026932593,-0.0114784240785799,0.0693714855540307),
    @(0.0941223041295222,0.192308607795969,-0.00168303603910866,0.0250475264455221),
    @(0.0791132385922413,0.179001779279827,-0.00448644789256815,0.0254817705835679),
    @(0.114317161007836,0.22516188384634,-0.00319420306225574,0.0282767150098225)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $values[$i][0]
    $ws.Cells.Item($row, 5).Value = $values[$i][1]
    $ws.Cells.Item($row, 6).Value = $values[$i][2]
    $ws.Cells.Item($row, 7).Value = $values[$i][3]
}
